$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Fill in the 22/5/2025 row (row 21) of the tracking table with Thursday's
# labeling progress numbers.
$ws.Range("E21").Value = 178
$ws.Range("F21").Value = 339
$ws.Range("G21").Value = 49
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 650
$ws.Range("J21").Value = "N/A"

# Update the view: scroll so column F is the left-most visible column,
# and move the active selection to I25.
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("I25").Select()
